$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 152
$ws.Range("I4").Value = 122.2
$ws.Range("J4").Value = 450
$ws.Range("K4").Value = 122.2
$ws.Range("L4").Value = 450
$ws.Range("M4").Value = -8.200000000000003
$ws.Range("N4").Value = -678
$ws.Range("H17").Value = 2584.4
$ws.Range("J17").Value = 2584.4
$ws.Range("L17").Value = 7753.200000000001
$ws.Range("N17").Value = -8089.200000000001
$ws.Range("H18").Value = 10000.667
$ws.Range("I18").Value = 5000
$ws.Range("K18").Value = 5000
$ws.Range("M18").Value = -4716
$ws.Range("H62").Value = 2373
$ws.Range("I62").Value = 2403.1428
$ws.Range("J62").Value = 2353.818
$ws.Range("K62").Value = 2403.1428
$ws.Range("L62").Value = 2353.818
$ws.Range("M62").Value = -1779.1428
$ws.Range("N62").Value = -3601.818
$ws.Range("H65").Value = 2373
$ws.Range("I65").Value = 2403.1428
$ws.Range("J65").Value = 2353.818
$ws.Range("K65").Value = 12015.714
$ws.Range("L65").Value = 11769.09
$ws.Range("M65").Value = -8895.714
$ws.Range("N65").Value = -18009.09
$ws.Range("H92").Value = 787
$ws.Range("I92").Value = 787
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 787
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 461
$ws.Range("N92").ClearContents()
$ws.Range("H138").Value = 2352.9285
$ws.Range("I138").Value = 882.3333
$ws.Range("K138").Value = 2646.9999
$ws.Range("M138").Value = 2493.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 874.5
$ws.Range("I4").Value = 866
$ws.Range("J4").Value = 900
$ws.Range("K4").Value = 866
$ws.Range("L4").Value = 900
$ws.Range("M4").Value = -750
$ws.Range("N4").Value = -1132
$ws.Range("H61").Value = 7011
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 5197.4614
$ws.Range("I74").Value = 5031.7
$ws.Range("K74").Value = 5031.7
$ws.Range("M74").Value = -4157.7
$ws.Range("H77").Value = 5197.4614
$ws.Range("I77").Value = 5031.7
$ws.Range("K77").Value = 25158.5
$ws.Range("M77").Value = -20790.5
$ws.Range("H94").Value = 87664.5
$ws.Range("J94").Value = 87664.5
$ws.Range("L94").Value = 87664.5
$ws.Range("N94").Value = -89466.5
$ws.Range("H110").Value = 3498.25
$ws.Range("I110").Value = 3666
$ws.Range("J110").Value = 2995
$ws.Range("K110").Value = 3666
$ws.Range("L110").Value = 2995
$ws.Range("M110").Value = -1621
$ws.Range("N110").Value = -7085
$ws.Range("H132").Value = 2171.6667
$ws.Range("I132").Value = 1750.5
$ws.Range("J132").Value = 3014
$ws.Range("K132").Value = 5251.5
$ws.Range("L132").Value = 9042
$ws.Range("M132").Value = -2721.5
$ws.Range("N132").Value = -14102
$ws.Range("H136").Value = 7011
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H22").Value = 361.33334
$ws.Range("I22").Value = 361.33334
$ws.Range("K22").Value = 361.33334
$ws.Range("M22").Value = -188.33334
$ws.Range("H86").Value = 2583.077
$ws.Range("J86").Value = 2569.1
$ws.Range("L86").Value = 2569.1
$ws.Range("N86").Value = -4815.1
$ws.Range("H89").Value = 2583.077
$ws.Range("J89").Value = 2569.1
$ws.Range("L89").Value = 12845.5
$ws.Range("N89").Value = -24077.5
$ws.Range("H95").Value = 21787
$ws.Range("J95").Value = 21787
$ws.Range("L95").Value = 21787
$ws.Range("N95").Value = -27279
$ws.Range("H134").Value = 11054.923
$ws.Range("I134").Value = 9476.166999999999
$ws.Range("K134").Value = 28428.501
$ws.Range("M134").Value = -25893.501

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4935.048
$ws.Range("I31").Value = 3948.3572
$ws.Range("K31").Value = 3948.3572
$ws.Range("M31").Value = -3653.3572
$ws.Range("H34").Value = 4935.048
$ws.Range("I34").Value = 3948.3572
$ws.Range("K34").Value = 3948.3572
$ws.Range("M34").Value = -3746.3572
$ws.Range("H41").Value = 19600
$ws.Range("J41").Value = 21666.666
$ws.Range("L41").Value = 21666.666
$ws.Range("N41").Value = -22522.666
$ws.Range("H56").Value = 52000
$ws.Range("J56").Value = 52000
$ws.Range("L56").Value = 52000
$ws.Range("N56").Value = -53690
$ws.Range("H59").Value = 29438.889
$ws.Range("I59").Value = 22487.5
$ws.Range("K59").Value = 22487.5
$ws.Range("M59").Value = -21342.5
$ws.Range("H60").Value = 20800
$ws.Range("H62").Value = 5699.75
$ws.Range("I62").Value = 5449.5
$ws.Range("J62").Value = 5950
$ws.Range("K62").Value = 5449.5
$ws.Range("L62").Value = 5950
$ws.Range("M62").Value = -4825.5
$ws.Range("N62").Value = -7198
$ws.Range("H65").Value = 5699.75
$ws.Range("I65").Value = 5449.5
$ws.Range("J65").Value = 5950
$ws.Range("K65").Value = 27247.5
$ws.Range("L65").Value = 29750
$ws.Range("M65").Value = -24127.5
$ws.Range("N65").Value = -35990
$ws.Range("H68").Value = 27287.062
$ws.Range("I68").Value = 17399.223
$ws.Range("K68").Value = 17399.223
$ws.Range("M68").Value = -16650.223
$ws.Range("H71").Value = 27287.062
$ws.Range("I71").Value = 17399.223
$ws.Range("K71").Value = 52197.66900000001
$ws.Range("M71").Value = -48453.66900000001
$ws.Range("H93").Value = 19000
$ws.Range("I93").Value = 19000
$ws.Range("K93").Value = 19000
$ws.Range("M93").Value = -17128
$ws.Range("H107").Value = 941.125
$ws.Range("I107").Value = 1061.2858
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 1061.2858
$ws.Range("L107").Value = 100
$ws.Range("M107").Value = 858.7141999999999
$ws.Range("N107").Value = -3940
$ws.Range("H132").Value = 1282.5714
$ws.Range("I132").Value = 1010.17645
$ws.Range("J132").Value = 2440.25
$ws.Range("K132").Value = 3030.52935
$ws.Range("L132").Value = 7320.75
$ws.Range("M132").Value = -500.5293500000002
$ws.Range("N132").Value = -12380.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 428.4
$ws.Range("I36").Value = 455.5
$ws.Range("K36").Value = 1366.5
$ws.Range("M36").Value = -1197.5
$ws.Range("H75").Value = 4407
$ws.Range("J75").Value = 4407
$ws.Range("L75").Value = 13221
$ws.Range("N75").Value = -15217
$ws.Range("H78").Value = 4407
$ws.Range("J78").Value = 4407
$ws.Range("L78").Value = 39663
$ws.Range("N78").Value = -49647
$ws.Range("H80").Value = 9119.4
$ws.Range("I80").Value = 1798
$ws.Range("K80").Value = 5394
$ws.Range("M80").Value = -4458
$ws.Range("H83").Value = 9119.4
$ws.Range("I83").Value = 1798
$ws.Range("K83").Value = 16182
$ws.Range("M83").Value = -11502
$ws.Range("H109").Value = 1106.875
$ws.Range("I109").Value = 1106.875
$ws.Range("K109").Value = 3320.625
$ws.Range("M109").Value = -2280.625
$ws.Range("H114").Value = 981.5
$ws.Range("I114").Value = 619.3333
$ws.Range("J114").Value = 1343.6666
$ws.Range("K114").Value = 1857.9999
$ws.Range("L114").Value = 4030.9998
$ws.Range("M114").Value = 1396.0001
$ws.Range("N114").Value = -10538.9998
$ws.Range("H117").Value = 470.08334
$ws.Range("I117").Value = 343.25
$ws.Range("J117").Value = 723.75
$ws.Range("K117").Value = 1029.75
$ws.Range("L117").Value = 2171.25
$ws.Range("M117").Value = 2412.25
$ws.Range("N117").Value = -9055.25
$ws.Range("H121").Value = 1011.875
$ws.Range("I121").Value = 865.8
$ws.Range("J121").Value = 1078.2727
$ws.Range("K121").Value = 2597.4
$ws.Range("L121").Value = 3234.8181
$ws.Range("M121").Value = -1287.4
$ws.Range("N121").Value = -5854.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 148.25
$ws.Range("J107").Value = 399
$ws.Range("L107").Value = 399
$ws.Range("N107").Value = -4239
$ws.Range("H132").Value = 2333.5
$ws.Range("I132").Value = 2333.5
$ws.Range("K132").Value = 7000.5
$ws.Range("M132").Value = -4470.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2846.9
$ws.Range("I22").Value = 2916.6667
$ws.Range("J22").Value = 2742.25
$ws.Range("K22").Value = 2916.6667
$ws.Range("L22").Value = 2742.25
$ws.Range("M22").Value = -2621.6667
$ws.Range("N22").Value = -3332.25
$ws.Range("H27").Value = 2846.9
$ws.Range("I27").Value = 2916.6667
$ws.Range("J27").Value = 2742.25
$ws.Range("K27").Value = 2916.6667
$ws.Range("L27").Value = 2742.25
$ws.Range("M27").Value = -2809.6667
$ws.Range("N27").Value = -2956.25
$ws.Range("H136").Value = 3502
$ws.Range("I136").Value = 3502
$ws.Range("K136").Value = 10506
$ws.Range("M136").Value = -7956

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 9997.5
$ws.Range("J94").Value = 9997.5
$ws.Range("L94").Value = 9997.5
$ws.Range("N94").Value = -11799.5
$ws.Range("H136").Value = 1352.25
$ws.Range("I136").Value = 952
$ws.Range("J136").Value = 1752.5
$ws.Range("K136").Value = 2856
$ws.Range("L136").Value = 5257.5
$ws.Range("M136").Value = -306
$ws.Range("N136").Value = -10357.5
